$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Sheet "ip_address_list": rows 9-12 rotate (old row 9 -> new row 12; old
# rows 10/11/12 shift up to 9/10/11), plus row 14 notes get trimmed.
# -------------------------------------------------------------------------
$wsIp = $wb.Worksheets.Item("ip_address_list")

# New row 9 (was old row 10: 511_Teleflex)
$wsIp.Range("A9").Value = "511_Teleflex"
$wsIp.Range("B9").Value = "192.168.1.242"
$wsIp.Range("D9").Value = "Teleflex "

# New row 10 (was old row 11: 503_Witte)
$wsIp.Range("A10").Value = "503_Witte"
$wsIp.Range("B10").Value = "192.168.0.240"
$wsIp.Range("D10").Value = "PC:`t10.96.205.175`nNAS:`t10.96.205.166`nFH:`t10.96.205.154`n`t10.96.205.267`n-----------------------------------------`nuser:JHV_Vision, omron `nPass:*Jhv2708`n---------------------------------------`nFortiClient Austin: `nPass:`n1Pm#J@PFIkzM&Q@i `nUVt1@Ex2p78kxp30atD7we@!qGK"

# New row 11 (was old row 12: 497_Edcha)
$wsIp.Range("A11").Value = "497_Edcha"
$wsIp.Range("B11").Value = "172.26.7.240"
$wsIp.Range("D11").Value = "FortiClient Edcha Ex2p78kxp30"
$wsIp.Range("E11").Value = 0

# New row 12 (was old row 9: 514_Teleflex)
$wsIp.Range("A12").Value = "514_Teleflex"
$wsIp.Range("B12").Value = "192.168.14.240"
$wsIp.Range("D12").Value = "PC:192.168.14.240`nCAM: 192.168.14.??NAS:192.168.14.245`n*******************************`nuser: Vision`npass: *Jhv2708"
$wsIp.Range("E12").Value = 1

# Row 14: trim the FortiClient Austin password lines out of the note.
$wsIp.Range("D14").Value = "FortiClient Austin: `nFH-2050-20`n10.96.205.80"

# -------------------------------------------------------------------------
# Sheet "disk_list": drop the "xfdx" test row (row 1) entirely - remaining
# rows shift up, and the hyperlink that lived on old C1 goes away with it.
# -------------------------------------------------------------------------
$wsDisk = $wb.Worksheets.Item("disk_list")
$wsDisk.Range("C1").Hyperlinks.Delete()
$wsDisk.Rows.Item(1).Delete()

# -------------------------------------------------------------------------
# Sheet "Settings": swap the two startup-window flags.
# -------------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("B3").Value = 1
$wsSettings.Range("B4").Value = 0

# -------------------------------------------------------------------------
# Sheet "projects_bin2" (hidden): the "xfdx" test row reappears here as the
# new row 2, between the existing rows 1 and 3.
# -------------------------------------------------------------------------
$wsBin = $wb.Worksheets.Item("projects_bin2")
$wsBin.Range("A2").Value = "xfdx"
$wsBin.Range("B2").Value = "P"
$wsBin.Range("C2").Value = "\\192.168.000.000\"
$wsBin.Range("D2").Value = "ss"
